# Update daily COVID stats for Slovakia (commit: "Updated: ut 26. 01. 2021")
# Applies new values to columns H (AgTests) and I (AgPosit) for a range
# of rows, matching the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 289; H = 64395;  I = 3699  },
    @{ Row = 292; H = 81241;  I = 7191  },
    @{ Row = 299; H = 64055;  I = 6713  },
    @{ Row = 306; H = 70729;  I = 7182  },
    @{ Row = 309; H = 57315;  I = 3962  },
    @{ Row = 310; H = 91076;  I = 5190  },
    @{ Row = 311; H = 37117;  I = 1395  },
    @{ Row = 313; H = 72960;  I = 3553  },
    @{ Row = 314; H = 65084;  I = 3346  },
    @{ Row = 315; H = 66203;  I = 3011  },
    @{ Row = 316; H = 49037;  I = 2278  },
    @{ Row = 317; H = 61174;  I = 2124  },
    @{ Row = 318; H = 25032;  I = 925   },
    @{ Row = 320; H = 86146;  I = 3892  },
    @{ Row = 321; H = 88965;  I = 2766  },
    @{ Row = 322; H = 104338; I = 2288  },
    @{ Row = 323; H = 147242; I = 2294  },
    @{ Row = 324; H = 221337; I = 2590  },
    @{ Row = 325; H = 593430; I = 5572  },
    @{ Row = 326; H = 326887; I = 2687  }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 8).Value = $u.H
    $ws.Cells.Item($u.Row, 9).Value = $u.I
}
